$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "68.809.15"
$ws.Range("E2").Value = "  -0.49%  "
$ws.Range("D3").Value = "3.864.85"
$ws.Range("E3").Value = "  +3.03%  "
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("D5").Value = "'600.41"
$ws.Range("E5").Value = "  -0.27%  "
$ws.Range("D6").Value = "'162.44"
$ws.Range("E6").Value = "  -2.78%  "
$ws.Range("D7").Value = "3.861.60"
$ws.Range("E7").Value = "  +2.97%  "
$ws.Range("E8").Value = "  +0.02%  "
$ws.Range("E9").Value = "  -1.74%  "
$ws.Range("E10").Value = "  -0.76%  "
$ws.Range("E11").Value = "  -1.21%  "
$ws.Range("E12").Value = "  -0.47%  "
$ws.Range("D13").Value = "'36.88"
$ws.Range("E13").Value = "  -2.87%  "
$ws.Range("E14").Value = "  -2.03%  "
$ws.Range("D15").Value = "4.515.51"
$ws.Range("E15").Value = "  +3.16%  "
$ws.Range("D16").Value = "3.867.55"
$ws.Range("E16").Value = "  +3.25%  "
$ws.Range("D17").Value = "68.987.79"
$ws.Range("E17").Value = "  -0.24%  "
$ws.Range("E18").Value = "  +2.58%  "
$ws.Range("E19").Value = "  -0.47%  "
$ws.Range("D20").Value = "'11.37"
$ws.Range("E20").Value = "  +2.72%  "
$ws.Range("D21").Value = "'17.09"
$ws.Range("E21").Value = "  -1.67%  "
$ws.Range("D22").Value = "'483.61"
$ws.Range("E22").Value = "  -1.93%  "
$ws.Range("D23").Value = "'0.718"
$ws.Range("E23").Value = "  -1.29%  "
$ws.Range("D24").Value = "'0.0000161"
$ws.Range("E24").Value = "  +6.61%  "
$ws.Range("D25").Value = "'83.95"
$ws.Range("E25").Value = "  -1.10%  "
$ws.Range("E26").Value = "  -2.87%  "
$ws.Range("D27").Value = "'12.07"
$ws.Range("E28").Value = "  -0.08%  "
$ws.Range("E29").Value = "  -1.53%  "
$ws.Range("E30").Value = "  -1.14%  "
$ws.Range("D31").Value = "4.018.40"
$ws.Range("E31").Value = "  +3.09%  "
$ws.Range("D32").Value = "'7.87"
$ws.Range("E32").Value = "  -3.08%  "
$ws.Range("D33").Value = "'32.25"
$ws.Range("E33").Value = "  +2.33%  "
$ws.Range("D34").Value = "'2.36"
$ws.Range("E34").Value = "  -4.20%  "
$ws.Range("D35").Value = "3.812.86"
$ws.Range("E35").Value = "  +3.51%  "
$ws.Range("D36").Value = "'0.106"
$ws.Range("E36").Value = "  -1.71%  "
$ws.Range("E37").Value = "  +1.90%  "
$ws.Range("E38").Value = "  +1.60%  "
$ws.Range("D39").Value = "'5.88"
$ws.Range("E39").Value = "  -1.56%  "
$ws.Range("E41").Value = "  -2.59%  "
$ws.Range("E42").Value = "  -1.93%  "
$ws.Range("D43").Value = "'432.12"
$ws.Range("E43").Value = "  +1.69%  "
$ws.Range("E45").Value = "  +0.12%  "
$ws.Range("D47").Value = "'8.38"
$ws.Range("E47").Value = "  -1.07%  "
$ws.Range("D48").Value = "'143.38"
$ws.Range("E48").Value = "  +1.32%  "
$ws.Range("D49").Value = "2.838.78"
$ws.Range("B50").Value = "VeChain"
$ws.Range("C50").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D50").Value = "'0.0357"
$ws.Range("E50").Value = "  +1.03%  "
$ws.Range("B51").Value = "EnergySwap"
$ws.Range("C51").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D51").Value = "'25.92"
$ws.Range("E51").Value = "  +12.55%  "
